# The workbook originally has, in tab order:
#   MAIN PAGE, START, find and replace, formatting, auto filter,
#   text to columns, home and end keys, remove duplicates,
#   paste formatting, paste formatting 2, Bold Italics Underline, FINAL
#
# Target layout:
#   MAIN PAGE, START, find and replace, formatting 1, formatting 2,
#   auto filter, text to columns, home and end keys, remove duplicates,
#   paste formatting, paste formatting 2, FINAL
#
# i.e. "Bold Italics Underline" is renamed to "formatting 1" and moved to
# sit right after "find and replace" (and right before the renamed
# "formatting" -> "formatting 2" sheet).

$wb = $excel.ActiveWorkbook

# Rename "Bold Italics Underline" -> "formatting 1" and relocate it so it
# becomes the sheet immediately after "find and replace".
$formatting1 = $wb.Worksheets.Item("Bold Italics Underline")
$formatting1.Name = "formatting 1"
$formatting1.Move($wb.Worksheets.Item("find and replace").Next())

# Rename the old "formatting" sheet -> "formatting 2" (it now immediately
# follows "formatting 1").
$formatting2 = $wb.Worksheets.Item("formatting")
$formatting2.Name = "formatting 2"
